# Making image acquisition_datetime optional
#
# This adds two new header columns - "validated" and "validation_datetime" -
# immediately before the trailing "comment" column on each of the four
# "*Output" dataset sheets (FieldIlluminationOutput, PSFBeadsOutput,
# ArgolightBOutput, ArgolightEOutput). The existing "comment" column is
# shifted two columns to the right to make room.

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "FieldIlluminationOutput",
    "PSFBeadsOutput",
    "ArgolightBOutput",
    "ArgolightEOutput"
)

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Locate the last used column in row 1 (this holds the "comment" header).
    $lastCol = $ws.Cells.Item(1, $ws.UsedRange.Columns.Count)
    $lastColIndex = $lastCol.Column

    # Insert two new blank columns right before the "comment" column,
    # shifting "comment" two columns to the right.
    $insertRange = $ws.Range($ws.Cells.Item(1, $lastColIndex), $ws.Cells.Item(1, $lastColIndex + 1))
    $insertRange.EntireColumn.Insert()

    # Fill in the two new header cells.
    $ws.Cells.Item(1, $lastColIndex).Value = "validated"
    $ws.Cells.Item(1, $lastColIndex + 1).Value = "validation_datetime"
}
